$d = $word.ActiveDocument

# --- Change 1: insert a new "Meta description" paragraph right after the title paragraph ---
$titlePara = $d.Paragraphs.Item(1)
if ($titlePara.Range.Text.TrimEnd("`r", "`a") -ne "Play Bomber Fruit for Free - Review of Capecod's Fruit-Themed Slot") {
    throw "Unexpected content for title paragraph: " + $titlePara.Range.Text
}
$insertionPoint = $d.Range($titlePara.Range.End - 1, $titlePara.Range.End - 1)

$metaRestText = ": Read our review of Bomber Fruit from Capecod, a fruit-themed slot with good variety and special symbols. Play for free and win big!"
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$metaParaXml = '<w:p ' + $ns + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">' + $metaRestText + '</w:t></w:r></w:p>'
# A trailing empty paragraph is required so InsertXML creates a genuine new
# paragraph break instead of merging our runs into the following paragraph.
$trailerParaXml = '<w:p ' + $ns + '></w:p>'
$insertionPoint.InsertXML($metaParaXml + $trailerParaXml)

# The trailing placeholder paragraph from above is now paragraph 3; remove it
# again so the document falls back to its original paragraph after our new one.
$emptyTrailerPara = $d.Paragraphs.Item(3)
$emptyTrailerPara.Range.Delete()

# --- Change 2 & 3: the two paragraphs at the very end of the document ---
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$metaDescPara = $d.Paragraphs.Item($count)

# Sanity-check before mutating, so failures are obvious rather than silently wrong.
if ($dupTitlePara.Range.Text.TrimEnd("`r", "`a") -ne "Play Bomber Fruit for Free - Review of Capecod's Fruit-Themed Slot") {
    throw "Unexpected content for duplicate title paragraph: " + $dupTitlePara.Range.Text
}
if ($metaDescPara.Range.Text.TrimEnd("`r", "`a") -ne "Read our review of Bomber Fruit from Capecod, a fruit-themed slot with good variety and special symbols. Play for free and win big!") {
    throw "Unexpected content for meta-description paragraph: " + $metaDescPara.Range.Text
}

# Change 2: delete the duplicated bold title paragraph near the bottom entirely.
$dupTitlePara.Range.Delete()

# Change 3: replace the meta-description paragraph's text with the DALLE image prompt,
# keeping the run's existing (italic) formatting intact.
$count = $d.Paragraphs.Count
$metaDescPara = $d.Paragraphs.Item($count)
$newPromptText = 'Prompt: DALLE, create a cartoon-style image featuring a happy Maya warrior with glasses for the game "Bomber Fruit". The image should include the Maya warrior holding a bomb and a basket of fruits in his other hand. The background should be colorful and incorporate some of the fruits from the game such as watermelon, lemon, and cherry. Make sure the image is eye-catching and playful, while still capturing the essence of the game.'
$textRange = $d.Range($metaDescPara.Range.Start, $metaDescPara.Range.End - 1)
$textRange.Text = $newPromptText
